$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as exact text (avoids Excel auto-converting
# numeric-looking strings like "1.001" or "29.355.68" into numbers),
# then restore the default "Normal" style so no stray number-format
# styling is left behind on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '29.355.68'
$ws.Range("E2").Value = '  +0.18%  '

Set-TextValue $ws.Range("D3") '1.879.35'
$ws.Range("E3").Value = '  +0.23%  '

Set-TextValue $ws.Range("D4") '1.001'
$ws.Range("E4").Value = '  +0.06%  '

Set-TextValue $ws.Range("D5") '0.7116'
$ws.Range("E5").Value = '  -0.12%  '

Set-TextValue $ws.Range("D6") '242.22'
$ws.Range("E6").Value = '  -0.08%  '

$ws.Range("E7").Value = '  +0.03%  '

Set-TextValue $ws.Range("D8") '0.08039'
$ws.Range("E8").Value = '  +3.97%  '

Set-TextValue $ws.Range("D9") '0.3130'
$ws.Range("E9").Value = '  +0.62%  '

Set-TextValue $ws.Range("D10") '25.25'
$ws.Range("E10").Value = '  +0.59%  '

Set-TextValue $ws.Range("D11") '0.08370'
$ws.Range("E11").Value = '  -0.87%  '

Set-TextValue $ws.Range("D12") '1.870.32'
$ws.Range("E12").Value = '  -0.42%  '

Set-TextValue $ws.Range("D13") '5.257'
$ws.Range("E13").Value = '  +0.82%  '

Set-TextValue $ws.Range("D14") '0.7197'
$ws.Range("E14").Value = '  +1.25%  '

Set-TextValue $ws.Range("D16") '6.256'
$ws.Range("E16").Value = '  +4.50%  '

Set-TextValue $ws.Range("D17") '0.000008370'
$ws.Range("E17").Value = '  +0.97%  '

Set-TextValue $ws.Range("D18") '29.352.18'

Set-TextValue $ws.Range("D19") '240.93'
$ws.Range("E19").Value = '  -0.74%  '

$ws.Range("E20").Value = '  +0.30%  '

Set-TextValue $ws.Range("D21") '2.127.54'
$ws.Range("E21").Value = '  -0.48%  '

Set-TextValue $ws.Range("D22") '1.001'
$ws.Range("E22").Value = '  +0.09%  '

Set-TextValue $ws.Range("D23") '7.804'
$ws.Range("E23").Value = '  -0.09%  '

$ws.Range("E24").Value = '  +0.10%  '

Set-TextValue $ws.Range("D25") '0.1591'
$ws.Range("E25").Value = '  -1.75%  '

Set-TextValue $ws.Range("D26") '163.01'
$ws.Range("E26").Value = '  -0.12%  '

Set-TextValue $ws.Range("D27") '9.057'
$ws.Range("E27").Value = '  +0.42%  '

$ws.Range("E28").Value = '  +0.28%  '

$ws.Range("E29").Value = '  -0.01%  '

Set-TextValue $ws.Range("D30") '4.422'
$ws.Range("E30").Value = '  +0.05%  '

Set-TextValue $ws.Range("D31") '4.339'
$ws.Range("E31").Value = '  +0.23%  '

Set-TextValue $ws.Range("D32") '1.189'
$ws.Range("E32").Value = '  -7.56%  '

Set-TextValue $ws.Range("D33") '0.05377'
$ws.Range("E33").Value = '  +2.28%  '

Set-TextValue $ws.Range("D34") '1.949'

$ws.Range("E35").Value = '  +0.55%  '

Set-TextValue $ws.Range("D36") '0.7512'
$ws.Range("E36").Value = '  +0.61%  '

Set-TextValue $ws.Range("D37") '2.696'
$ws.Range("E37").Value = '  +0.40%  '

Set-TextValue $ws.Range("D38") '1.296.01'
$ws.Range("E38").Value = '  +11.58%  '

Set-TextValue $ws.Range("D39") '0.01886'
$ws.Range("E39").Value = '  +1.46%  '

Set-TextValue $ws.Range("D40") '2.739'
$ws.Range("E40").Value = '  +0.75%  '

Set-TextValue $ws.Range("D41") '6.581'
$ws.Range("E41").Value = '  +3.34%  '

Set-TextValue $ws.Range("D42") '110.69'
$ws.Range("E42").Value = '  +3.95%  '

$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D43") '73.25'
$ws.Range("E43").Value = '  +0.25%  '

$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D44") '0.8911'
$ws.Range("E44").Value = '  +0.19%  '

$ws.Range("E45").Value = '  +8.55%  '

$ws.Range("E46").Value = '  +0.08%  '

Set-TextValue $ws.Range("D47") '2.025.74'
$ws.Range("E47").Value = '  +0.06%  '

$ws.Range("E48").Value = '  -0.26%  '

Set-TextValue $ws.Range("D49") '0.5210'
$ws.Range("E49").Value = '  +0.32%  '

Set-TextValue $ws.Range("D50") '9.471'
$ws.Range("E50").Value = '  +0.82%  '

Set-TextValue $ws.Range("D51") '0.4363'
$ws.Range("E51").Value = '  +1.53%  '
